$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The roster rows for "Zeke Nnaji" (row 10) and "Jeff Green" (row 11) were
# swapped (everything except the "No." column A, which stays sequential).

$cols = @("B","C","D","E","F","G","H","J","K")
foreach ($col in $cols) {
    $addr10 = $col + "10"
    $addr11 = $col + "11"
    $v10 = $ws.Range($addr10).Value2
    $v11 = $ws.Range($addr11).Value2
    $ws.Range($addr10).Value = $v11
    $ws.Range($addr11).Value = $v10
}

# Column I ("Exp") holds text that can look numeric ("2", "14", "R", ...),
# so force text formatting before writing to avoid Excel auto-converting
# the swapped values into numeric cells.
$i10 = $ws.Range("I10").Value2
$i11 = $ws.Range("I11").Value2
$ws.Range("I10").NumberFormat = "@"
$ws.Range("I11").NumberFormat = "@"
$ws.Range("I10").Value = $i11
$ws.Range("I11").Value = $i10
$ws.Range("I10:I11").NumberFormat = "General"
$ws.Range("I10:I11").Style = "Normal"
